# Add four new login rows (8-11) plus one extra styled-but-empty row (12)
# to the Username/Password table on Sheet1, mirroring the existing table's
# formatting (style "4" used by A5:A7, style "2" used by A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new user "hhhhhhhhhh" / secret_sauce, formatted like A5:A7 ---
[void]$ws.Range("A7").Copy()
[void]$ws.Range("A8").PasteSpecial(-4122)    # xlPasteFormats - copy formatting only
$ws.Range("A8").Value = "hhhhhhhhhh"
$ws.Range("B8").Value = "secret_sauce"
$ws.Rows.Item(8).RowHeight = 15.5            # matches the taller rows above it

# --- Row 9: password-only row (no username) ---
$ws.Range("B9").Value = "secret_sauce"

# --- Row 10: standard_user, formatted like A2 ---
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "standard_user"

# --- Row 11: standard_user / new password "hhhhhh" ---
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "standard_user"
$ws.Range("B11").Value = "hhhhhh"

# --- Row 12: empty username cell but carries A2's formatting ---
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A12").PasteSpecial(-4122)

# Leave the final selection on A12, matching the saved workbook state.
[void]$ws.Range("A12").Select()
